{"js": "// Insert a new, centered \"INTRODUCCI\u00d3N\" heading paragraph before the\n// document's existing (first) paragraph, and center that existing\n// paragraph (which holds the _GoBack bookmark) as well.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Insert the new \"INTRODUCCI\u00d3N\" paragraph immediately before the first\n// (existing) paragraph.\nconst introParagraph = firstParagraph.insertParagraph(\"INTRODUCCI\u00d3N\", Word.InsertLocation.before);\nintroParagraph.alignment = Word.Alignment.centered;\n\n// Center the original (pre-existing) paragraph as well.\nfirstParagraph.alignment = Word.Alignment.centered;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Insert a new centered \"INTRODUCCI\u00d3N\" paragraph before the document's\n# existing (first) paragraph.\n$firstPara = $d.Paragraphs.First\n$firstRange = $firstPara.Range\n$firstRange.InsertParagraphBefore()\n$firstRange.InsertBefore(\"INTRODUCCI\u00d3N\")\n\n# The new heading paragraph is now the first paragraph; center it.\n$wdAlignParagraphCenter = 1\n$headingPara = $d.Paragraphs(1)\n$headingPara.Alignment = $wdAlignParagraphCenter\n\n# Center the original paragraph (now the second one, holding the bookmark).\n$bookmarkPara = $d.Paragraphs(2)\n$bookmarkPara.Alignment = $wdAlignParagraphCenter\n"}
